# Apply Simplified -> Traditional Chinese translation edits
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "英语" "英文"
Replace-Text " / 葡萄牙语 / 法语 / 泰语 / 越南语 / 西班牙语" " / 葡萄牙文 / 法文 / 泰文 / 越南文 / 西班牙文"
Replace-Text "简要" "簡介"
Replace-Text "发送给确认参加活动的与会者的电子邮件。 我们希望与他们分享航班和住宿的预订细节。" "發送給確認參加活動的參與者的電子郵件。 我們想與他們分享航班和住宿的預定詳情。"
Replace-Text "目标受众" "目標受眾"
Replace-Text "活动与会者" "活動參加者"
Replace-Text "主题: " "主題: "
Replace-Text "这是您关于" "這裡是您的"
Replace-Text "[活动名称]" "[活動名稱]"
Replace-Text "的预订详情 " " 預訂詳情 "
Replace-Text "我们迫不及待想与您会面！ " "我們迫不及待想見到您！ "
Replace-Text "[合作伙伴姓名]" "[合作夥伴姓名]"
Replace-Text ", " "， "
Replace-Text "希望您和我们一样对 " "希望您和我們一樣對 "
Replace-Text " 感到兴奋。 活动即将开始，我们已经做好了一切准备，希望您能和我们一起参加这次 " " 感到興奮。 活動即將開始，我們已經做好了一切準備，希望您能和我們一起參加這次 "
Replace-Text "会议/研讨会/旅行" "會議/研討會/旅行"
Replace-Text "在这封电子邮件中，我们链接/附上了以下文件：" "在此電子郵件中，我們已連結/附上以下文件："
Replace-Text "您的来回机票" "您的來回機票"
Replace-Text "您的住宿预订详情" "您的住宿預訂詳情"
Replace-Text "您的签证信息 " "您的簽證資料 "
Replace-Text "（如适用）" "（如適用）"
Replace-Text "如果有任何疑问，请联系我们：" "如有任何疑問，請與我們聯繫："
Replace-Text "实时聊天" "即時聊天"
Replace-Text " |" " 或 "
Replace-Text " WhatsApp" "WhatsApp"
Replace-Text ". " "。 "
Replace-Text "如有任何疑问，请通过 " "如有任何疑問，請通過 "
Replace-Text "[电子邮件地址]" "[電子郵件地址]"
Replace-Text "[WHATSAPP 号码]" "[WHATSAPP 號碼]"
Replace-Text " (WhatsApp) 联系您的区域经理 " " (WhatsApp) 聯繫您的區域經理 "
Replace-Text "[姓名]" "[NAME]"
Replace-Text "在 " "在"
Replace-Text " 见！" "見！"

function Replace-CommentText($old, $new) {
    for ($i = 1; $i -le $d.Comments.Count; $i++) {
        $c = $d.Comments.Item($i)
        if ($c.Range.Text -eq $old) {
            $c.Range.Text = $new
        }
    }
}

Replace-CommentText "选择其一" "選擇其中之一"
Replace-CommentText "检查这些是否是包含的文件" "確認這些是否是包含的文件"
Replace-CommentText "选择一个" "選擇一個"
